# Refresh Leve market-price derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job profit sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business
$ws.Range("H51").Value = 4500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4500
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5468

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 1107.8235
$ws.Range("J112").Value = 1107.8235
$ws.Range("L112").Value = 3323.4705
$ws.Range("N112").Value = -5539.470499999999

# Row 117: A Greater Grimoire
$ws.Range("H117").Value = 48634
$ws.Range("J117").Value = 48634
$ws.Range("L117").Value = 48634
$ws.Range("N117").Value = -57812

# Row 121: Mindful Medicine
$ws.Range("H121").Value = 1010853.5
$ws.Range("I121").Value = 265
$ws.Range("J121").Value = 1111912.4
$ws.Range("K121").Value = 795
$ws.Range("L121").Value = 3335737.2
$ws.Range("M121").Value = 952
$ws.Range("N121").Value = -3339231.2

# Row 123: Nearly Bare
$ws.Range("H123").Value = 38078.934
$ws.Range("J123").Value = 38078.934
$ws.Range("L123").Value = 38078.934
$ws.Range("N123").Value = -47878.934

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2604.4695
$ws.Range("I137").Value = 980.8261
$ws.Range("J137").Value = 3102.3867
$ws.Range("K137").Value = 2942.4783
$ws.Range("L137").Value = 9307.160100000001
$ws.Range("M137").Value = -392.4782999999998
$ws.Range("N137").Value = -14407.1601

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2526.5217
$ws.Range("I138").Value = 2323.875
$ws.Range("J138").Value = 2634.6
$ws.Range("K138").Value = 6971.625
$ws.Range("L138").Value = 7903.799999999999
$ws.Range("M138").Value = -1831.625
$ws.Range("N138").Value = -18183.8

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 3331
$ws.Range("I141").Value = 2344.077
$ws.Range("J141").Value = 8463
$ws.Range("K141").Value = 7032.231000000001
$ws.Range("L141").Value = 25389
$ws.Range("M141").Value = -1852.231000000001
$ws.Range("N141").Value = -35749

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 2178.5715
$ws.Range("I2").Value = 2216.6667
$ws.Range("K2").Value = 2216.6667
$ws.Range("M2").Value = -2103.6667

# Row 44: Very Slow Array
$ws.Range("H44").Value = 37704
$ws.Range("J44").Value = 37704
$ws.Range("L44").Value = 37704
$ws.Range("N44").Value = -38680

# Row 55: Employee Retention
$ws.Range("H55").Value = 27509
$ws.Range("I55").Value = 3048
$ws.Range("J55").Value = 51970
$ws.Range("K55").Value = 3048
$ws.Range("L55").Value = 51970
$ws.Range("M55").Value = -2733
$ws.Range("N55").Value = -52600

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2979
$ws.Range("I63").Value = 2433.125
$ws.Range("J63").Value = 5162.5
$ws.Range("K63").Value = 2433.125
$ws.Range("L63").Value = 5162.5
$ws.Range("M63").Value = -1747.125
$ws.Range("N63").Value = -6534.5

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2979
$ws.Range("I66").Value = 2433.125
$ws.Range("J66").Value = 5162.5
$ws.Range("K66").Value = 12165.625
$ws.Range("L66").Value = 25812.5
$ws.Range("M66").Value = -8733.625
$ws.Range("N66").Value = -32676.5

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1584.9166
$ws.Range("I74").Value = 896.4
$ws.Range("J74").Value = 2076.7144
$ws.Range("K74").Value = 896.4
$ws.Range("L74").Value = 2076.7144
$ws.Range("M74").Value = -22.39999999999998
$ws.Range("N74").Value = -3824.7144

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1584.9166
$ws.Range("I77").Value = 896.4
$ws.Range("J77").Value = 2076.7144
$ws.Range("K77").Value = 4482
$ws.Range("L77").Value = 10383.572
$ws.Range("M77").Value = -114
$ws.Range("N77").Value = -19119.572

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1677.9032
$ws.Range("I110").Value = 1565.6666
$ws.Range("K110").Value = 1565.6666
$ws.Range("M110").Value = 479.3334

# Row 116: No Scope
$ws.Range("H116").Value = 2178.5715
$ws.Range("I116").Value = 2216.6667
$ws.Range("K116").Value = 2216.6667
$ws.Range("M116").Value = 77.33329999999978

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3420
$ws.Range("I122").Value = 6050
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 18150
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -15700
$ws.Range("N122").Value = -9899.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 2178.5715
$ws.Range("I3").Value = 2216.6667
$ws.Range("K3").Value = 2216.6667
$ws.Range("M3").Value = -2102.6667

# Row 35: Lancers' Creed
$ws.Range("H35").Value = 34345
$ws.Range("J35").Value = 34345
$ws.Range("L35").Value = 34345
$ws.Range("N35").Value = -34965

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 3240.875
$ws.Range("I82").Value = 3240.875
$ws.Range("K82").Value = 3240.875
$ws.Range("M82").Value = -2857.875

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 3240.875
$ws.Range("I85").Value = 3240.875
$ws.Range("K85").Value = 3240.875
$ws.Range("M85").Value = -1914.875

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2468.4167
$ws.Range("I134").Value = 2052.3928
$ws.Range("J134").Value = 3924.5
$ws.Range("K134").Value = 6157.178400000001
$ws.Range("L134").Value = 11773.5
$ws.Range("M134").Value = -3622.178400000001
$ws.Range("N134").Value = -16843.5

# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 33147.375
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 33147.375
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 33147.375
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -43507.375

$ws = $wb.Worksheets.Item("CRP")
# Row 93: Reeling for Rods
$ws.Range("H93").Value = 11326.75
$ws.Range("I93").Value = 2653.5
$ws.Range("K93").Value = 2653.5
$ws.Range("M93").Value = -781.5

# Row 100: Run Before They Walk
$ws.Range("H100").Value = 46000
$ws.Range("J100").Value = 46000
$ws.Range("L100").Value = 46000
$ws.Range("N100").Value = -48164

$ws = $wb.Worksheets.Item("CUL")
# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 2900.8
$ws.Range("I98").Value = 2900
$ws.Range("J98").Value = 2901.3333
$ws.Range("K98").Value = 8700
$ws.Range("L98").Value = 8703.999899999999
$ws.Range("M98").Value = -7202
$ws.Range("N98").Value = -11699.9999

# Row 133: Friends Are Food
$ws.Range("H133").Value = 5093
$ws.Range("I133").Value = 3898.5715
$ws.Range("K133").Value = 11695.7145
$ws.Range("M133").Value = -6635.7145

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 803511.5600000001
$ws.Range("I2").Value = 860894.5
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 860894.5
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -860781.5
$ws.Range("N2").Value = -376

# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 5051.3403
$ws.Range("I70").Value = 5116.0312
$ws.Range("J70").Value = 4913.3335
$ws.Range("K70").Value = 5116.0312
$ws.Range("L70").Value = 4913.3335
$ws.Range("M70").Value = -4846.0312
$ws.Range("N70").Value = -5453.3335

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 5051.3403
$ws.Range("I73").Value = 5116.0312
$ws.Range("J73").Value = 4913.3335
$ws.Range("K73").Value = 5116.0312
$ws.Range("L73").Value = 4913.3335
$ws.Range("M73").Value = -4180.0312
$ws.Range("N73").Value = -6785.3335

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 4003.5
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 4014
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 4014
$ws.Range("M102").Value = -2378
$ws.Range("N102").Value = -7258

# Row 104: Speak Softly and Carry a Metal Rod
$ws.Range("H104").Value = 42165.75
$ws.Range("J104").Value = 42165.75
$ws.Range("L104").Value = 42165.75
$ws.Range("N104").Value = -49153.75

# Row 120: A Beneficent Elegy
$ws.Range("H120").Value = 22317
$ws.Range("J120").Value = 22317
$ws.Range("L120").Value = 22317
$ws.Range("N120").Value = -31993

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1900
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5700
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3250
$ws.Range("N122").ClearContents()

# Row 123: Workplace Workout
$ws.Range("H123").Value = 14271.667
$ws.Range("J123").Value = 14271.667
$ws.Range("L123").Value = 14271.667
$ws.Range("N123").Value = -19171.667

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 6095.6924
$ws.Range("I126").Value = 12272.4
$ws.Range("J126").Value = 2235.25
$ws.Range("K126").Value = 36817.2
$ws.Range("L126").Value = 6705.75
$ws.Range("M126").Value = -34347.2
$ws.Range("N126").Value = -11645.75

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1375.9231
$ws.Range("I22").Value = 1366.5555
$ws.Range("J22").Value = 1397
$ws.Range("K22").Value = 1366.5555
$ws.Range("L22").Value = 1397
$ws.Range("M22").Value = -1071.5555
$ws.Range("N22").Value = -1987

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1375.9231
$ws.Range("I27").Value = 1366.5555
$ws.Range("J27").Value = 1397
$ws.Range("K27").Value = 1366.5555
$ws.Range("L27").Value = 1397
$ws.Range("M27").Value = -1259.5555
$ws.Range("N27").Value = -1611

# Row 121: A Shoe In
$ws.Range("H121").Value = 30984.5
$ws.Range("J121").Value = 30984.5
$ws.Range("L121").Value = 30984.5
$ws.Range("N121").Value = -34478.5

$ws = $wb.Worksheets.Item("WVR")
# Row 121: Healing Headwear
$ws.Range("H121").Value = 43392
$ws.Range("J121").Value = 43392
$ws.Range("L121").Value = 43392
$ws.Range("N121").Value = -46886
